$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0005654550879112819
$ws.Range("B3").Value = 0.0152987983553583
$ws.Range("C3").Value = 14
$ws.Range("B4").Value = 0.07207341979134174
$ws.Range("B5").Value = 0.04425037842602819
$ws.Range("C5").Value = 11
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 21
$ws.Range("B7").Value = 0.0003203212830724667
$ws.Range("C7").Value = 17
$ws.Range("B8").Value = 0.003147217299551031
$ws.Range("C8").Value = 16
$ws.Range("C9").Value = 23
$ws.Range("B10").Value = 0.00231072907423524
$ws.Range("C10").Value = 21
$ws.Range("B11").Value = 0.00207724329404968
$ws.Range("C11").Value = 17
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 13
$ws.Range("C13").Value = 23
$ws.Range("B14").Value = 0.0004063500479780657
$ws.Range("C14").Value = 21
$ws.Range("C15").Value = 4
$ws.Range("B16").Value = 0.02258330274355114
$ws.Range("B18").Value = 0.006922060751502812
$ws.Range("B19").Value = 0.00700913032577688
$ws.Range("C19").Value = 10
$ws.Range("B20").Value = 0.01043456322946934
$ws.Range("C20").Value = 6
$ws.Range("B21").Value = 0.001829928053239449
$ws.Range("B22").Value = 0.001297756900458857
$ws.Range("B23").Value = 0.02028643846604705
$ws.Range("B24").Value = 0.04050943981025679
$ws.Range("B25").Value = 0.03647524836086635
$ws.Range("B26").Value = 0.04684119439717893
$ws.Range("B27").Value = 0.04636743495132034
$ws.Range("B28").Value = 0.05639383736517622
$ws.Range("C28").Value = 3
$ws.Range("B29").Value = 0.05666897832244366
$ws.Range("B30").Value = 0.05261479791126966
$ws.Range("B31").Value = 0.06167280379676754
$ws.Range("B32").Value = 0.06226873842624325
$ws.Range("B33").Value = 0.008303656323646713
$ws.Range("B34").Value = 0.007256950239384885
$ws.Range("B35").Value = 0.01163556754276041
$ws.Range("B36").Value = 0.008569525125595057
$ws.Range("B37").Value = 0.01518035013359142
$ws.Range("B38").Value = 0.01345291701184426
